$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Row 2: update the existing data values
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 2083885
$ws.Range("B2").Value = 168918
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 3.5

# ---------------------------------------------------------------------------
# Rows 3-5: these were placeholder/blank rows - fill in the new data rows
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2089334
$ws.Range("B3").Value = 16392
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3.5

$ws.Range("A4").Value = 2093685
$ws.Range("B4").Value = 52719
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 3.5

$ws.Range("A5").Value = 2101450
$ws.Range("B5").Value = 58407
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 3.5

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------
# D3:D5 / E3:E5 should carry the same number format / font / border as the
# already-correctly-styled D2 / E2 cells.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3:D5").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3:E5").PasteSpecial($xlPasteFormats) | Out-Null

# B2 switches away from the explicit-black "from_cred" look to the regular
# (theme-coloured) font used elsewhere (e.g. B4/B5), and gains word-wrap.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B2").WrapText = $true

# B3 takes on that same new look (font + number format + wrap).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection moves from F11 to E10
# ---------------------------------------------------------------------------
$ws.Range("E10").Select() | Out-Null
